# Adds "SamplesTab" and "FilesTab" rows (3 and 4) to the "startup" sheet,
# mirroring the existing "CasesTab" row (2): same StatQuery / Neo4jData /
# WebData columns (C/D/E), new per-tab query text in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new tab names first (column A) so the shared-string table ends
# up with "SamplesTab" / "FilesTab" added before the two long query strings,
# matching the order new strings were introduced upstream.
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# Row 3 = SamplesTab
$ws.Range("B3").Value = "MATCH (ss:study_subject)`nWITH COLLECT(ss.study_subject_id) AS all_subjects`nMATCH (samp:sample)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nWHERE s.study_acronym IN [`"B`"]  `nWITH`n    distinct lp,`n    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,`n    collect(distinct f.file_id) AS files,`n    samp, ss, s, p, all_subjects`nRETURN`n samp.sample_id AS ``Sample ID``,`n            ss.study_subject_id AS ``Case ID``,`n            p.program_acronym AS ``Program Code``,`n            s.study_acronym AS ``Arm``,`n            ss.disease_subtype AS ``Diagnosis``,`n            samp.tissue_type AS ``Tissue Type``,`n            samp.composition AS ``Tissue Composition``,`n            samp.sample_anatomic_site AS ``Sample Anatomic Site``,`n            samp.method_of_sample_procurement AS ``Sample Procurement Method``"
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Row 4 = FilesTab
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nMATCH (f)-[:file_of_sample]->(samp)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)`nMATCH (s)-[:study_of_program]->(p)`nMATCH (d)-[:diagnosis_of_study_subject]->(ss)`nMATCH (tp)-[:tp_of_diagnosis]->(d)`nWHERE s.study_acronym IN [`"B`"] `n WITH`n        f, parent,p, ss, d,tp, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent,p, ss, d,tp, s, samp,`n        f.file_size /(1024^i) AS value,`n        10^precision AS factor,`n        units[i] as unit`nWITH`n        f, parent,p, ss, d,tp, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN Distinct`n    f.file_name AS ``File Name``,`n    head(labels(samp)) AS ``Association``,`n    f.file_description AS ``Description``,`n    f.file_format AS ``File Format``,`n     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    p.program_acronym AS ``Program Code``,`n    s.study_acronym AS ``Arm``,`n    ss.study_subject_id AS ``Case ID``,`n    samp.sample_id AS ``Sample ID```n    order by f.file_name"
$ws.Range("C4").Value = $ws.Range("C2").Value()
$ws.Range("D4").Value = $ws.Range("D2").Value()
$ws.Range("E4").Value = $ws.Range("E2").Value()
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Selection / zoom, as left by the edit.
$ws.Range("C4").Select()
$excel.ActiveWindow.Zoom = 60

# Row heights auto-fit by Excel for the wrapped, multi-line query cells.
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6
